# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on the "Rules" sheet switches from the text "R40" to
# the text "1" (still stored as text, not as a number).
#
# A plain  $ws.Range("B11").Value = "1"  would let Excel's type-inference
# treat "1" as a genuine number (losing the text semantics). Writing the
# digit through a scratch cell that already carries an explicit text
# number-format and then pasting *values only* into B11 keeps B11's
# original style/format untouched while still landing a literal text
# value in the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$target = $ws.Range("B11")

# Use a far-away, otherwise empty cell as scratch space so nothing else
# on the sheet is disturbed.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues - value/type only, keeps B11's own style

# Remove every trace of the scratch cell again.
$scratch.Clear()
